# Edit script: apply the changes described by the diff to 02_Basic_Operators.pptx
#  1. Slide 12, body placeholder: replace the single long paragraph with four
#     shorter paragraphs (new wording) and per-paragraph spacing.
#  2. Slide 9, table (graphicFrame shape 3): change the table style GUID.

$p = $ppt.ActivePresentation

# --- 1. Slide 12 body text -------------------------------------------------
$slide12 = $p.Slides.Item(12)
$body = $slide12.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

$para1 = "erstelle die variable a und initialisiere sie mit 5"
$para2 = "erstelle die variable b und initialisiere sie mit 10"
$para3 = "überprüfe ob a in der Range 1...10 ist und b nicht nil ist"
$para4 = "Wenn dies True ist gib a + b aus"

$tr.Text = $para1 + "`r" + $para2 + "`r" + $para3 + "`r" + $para4

# Fix up the per-paragraph spacing (SpaceBefore/SpaceAfter are in points;
# 16pt == spcPts val="1600" in the underlying OOXML).
$start = 1

$p1 = $tr.Characters($start, $para1.Length)
$p1.ParagraphFormat.SpaceBefore = 0
$p1.ParagraphFormat.SpaceAfter = 0
$start = $start + $para1.Length + 1

$p2 = $tr.Characters($start, $para2.Length)
$p2.ParagraphFormat.SpaceBefore = 16
$p2.ParagraphFormat.SpaceAfter = 0
$start = $start + $para2.Length + 1

$p3 = $tr.Characters($start, $para3.Length)
$p3.ParagraphFormat.SpaceBefore = 16
$p3.ParagraphFormat.SpaceAfter = 0
$start = $start + $para3.Length + 1

$p4 = $tr.Characters($start, $para4.Length)
$p4.ParagraphFormat.SpaceBefore = 16
$p4.ParagraphFormat.SpaceAfter = 16

# --- 2. Slide 9 table style --------------------------------------------------
$slide9 = $p.Slides.Item(9)
$tableShape = $slide9.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{AFA38E8C-6406-4A0C-AD43-D798152BFC0B}")
